# Assignment-1/Report.xlsx — "final submission with report downloaded from colab"
#
# The three sheets (Telco, Adult, Credit) hold model-evaluation numbers.
# This re-submission refreshes the measured accuracy / TPR / TNR / F1 /
# boosting-round numbers with a freshly re-run experiment, and nudges the
# Telco sheet's on-screen selection down one row (C17 -> C18) to match the
# other two sheets. (Workbook/style metadata churn visible in the XML diff
# - stripped xr:* revision GUIDs, reshuffled cellXfs entries that resolve
# to byte-identical visual formatting, calcId, fileVersion, etc. - is an
# artifact of the file having been regenerated by an external tool rather
# than an in-app formatting edit, so there is no corresponding user action
# to replay here.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Telco sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Telco")

$ws.Range("C4").Value  = 73.88
$ws.Range("C5").Value  = 51
$ws.Range("C6").Value  = 83
$ws.Range("C9").Value  = 0.52

$ws.Range("B14").Value = 78.43
$ws.Range("C14").Value = 77.57
$ws.Range("B15").Value = 78.51
$ws.Range("C15").Value = 78.57
$ws.Range("B16").Value = 79.18
$ws.Range("C16").Value = 79.06
$ws.Range("B17").Value = 79.64
$ws.Range("C17").Value = 79.49

# Move the saved selection/active cell on this sheet from C17 to C18
# (visible in the diff's <selection activeCell="C17".../> -> "C18").
$ws.Activate() | Out-Null
$ws.Range("C18").Select() | Out-Null

# ---------------------------------------------------------------------
# Adult sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Adult")

$ws.Range("B14").Value = 84.26
$ws.Range("C14").Value = 84.55
$ws.Range("B15").Value = 84.67
$ws.Range("C15").Value = 84.92
$ws.Range("B16").Value = 85.2
$ws.Range("C16").Value = 85.52
$ws.Range("B17").Value = 85.31
$ws.Range("C17").Value = 85.34

# ---------------------------------------------------------------------
# Credit sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Credit")

$ws.Range("C4").Value  = 99
$ws.Range("C5").Value  = 79
$ws.Range("C6").Value  = 100
$ws.Range("C9").Value  = 0.79

$ws.Range("B14").Value = 99.52
$ws.Range("C14").Value = 99.46
$ws.Range("B15").Value = 99.52
$ws.Range("C15").Value = 99.46
$ws.Range("B16").Value = 99.52
$ws.Range("C16").Value = 99.46
$ws.Range("B17").Value = 99.54
$ws.Range("C17").Value = 99.49

# Restore the workbook's active tab to "Credit" (activeTab stays 2 /
# tabSelected stays on Credit throughout the original diff - only the
# Telco sheet's own saved selection moved).
$ws.Activate() | Out-Null
